$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegistrationData")

# Header: "Contact no" -> "Contact"
$ws.Range("D1").Value = "Contact"

# Fix D1 selection will be handled later

# Copy formatting from column C (Email) into column D (Contact) so the same
# hyperlink-style/border formatting used for rows 2-4 and 6 carries over.
$ws.Range("C1:C17").Copy()
$ws.Range("D1:D17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-apply header style (row 1 uses a distinct bold header style) since the
# format copy above would have overwritten it with column C's row1 style
# (they are actually identical style already, but keep value intact)
$ws.Range("D1").Value = "Contact"

# Restore/Set the contact numbers per row (D5 becomes blank, D6 becomes filled)
$ws.Range("D2").Value = 8447520166
$ws.Range("D3").Value = 8447520166
$ws.Range("D4").Value = 8447520166
$ws.Range("D5").ClearContents()
$ws.Range("D6").Value = 8447520166

# Hyperlinks for the Contact numbers (phone-like values hyperlinked to the
# registrant's email, matching the pattern already used for C/E/F columns)
$ws.Hyperlinks.Add($ws.Range("D2:D4"), "mailto:lokesh403@xtivia.com", "", "", "lokesh403@xtivia.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:lokesh403@xtivia.com", "", "", "lokesh403@xtivia.com")

# Hyperlinks.Add() replaces the anchor cell's stored value with the display
# text, so restore the numeric contact number afterwards.
$ws.Range("D2").Value = 8447520166
$ws.Range("D6").Value = 8447520166

# Hyperlinks.Add() also re-styles only the anchor cell of the range with a
# fresh style entry; reapply the shared hyperlink-cell format (from the
# adjoining rows of the same range, which already carry the correct style)
# so D2 and D6 end up visually consistent with D3/D4/C6.
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column widths: C and D share the bestFit width, E keeps width 16
$ws.Columns.Item(4).ColumnWidth = 20.25

# Update selection to D1
$ws.Range("D1").Select()
